$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '245.73'
$ws.Range("D2").ClearFormats()

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.463'
$ws.Range("D4").ClearFormats()

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05619'
$ws.Range("D5").ClearFormats()

# Row 6
$ws.Range("B6").Value = 'KuCoinToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.461'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '5KuCoinTokenKCS'

# Row 7
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8046'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '6MXTokenMX'

# Row 8
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.040'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '7FTXTokenFTT'

# Row 9
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1420'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '8WazirXWRX'

# Row 10
$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07298'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '9MandalaExchangeTokenMDX'

# Row 11
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.03181'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '10LiechtensteinCryptoassetsExchangeLCX'

# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.02937'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '11BitrueCoinBTR'

# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09261'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '12BitMartTokenBMX'

# Row 14
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001664'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '13BitForexTokenBF'

# Row 15
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.216'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '14MCDexMCB'

# Row 16
$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.04740'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '15CoinExTokenCET'

# Row 17
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.006339'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '16TigerCashTCH'

# Row 18
$ws.Range("B18").Value = 'HotbitToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.005071'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '17HotbitTokenHTB'

# Row 19
$ws.Range("B19").Value = 'BitKan'
$ws.Range("C19").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.001056'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '18BitKanKAN'

# Row 20
$ws.Range("B20").Value = 'NitroEx'
$ws.Range("C20").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0001503'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '19NitroExNTX'

# Row 21
$ws.Range("B21").Value = 'LEO'
$ws.Range("C21").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '3.988'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '20LEOLEO'

# Row 22
$ws.Range("B22").Value = 'GateToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.382'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '21GateTokenGT'

# Row 23
$ws.Range("B23").Value = 'BTSEToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.129'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '22BTSETokenBTSE'

# Row 24
$ws.Range("B24").Value = 'One'
$ws.Range("C24").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.01163'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '23OneONEBestin24h'

# Row 26
$ws.Range("E26").Value = '25ProBitTokenPROB'

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0002906'
$ws.Range("D27").ClearFormats()

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04158'
$ws.Range("D40").ClearFormats()

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006865'
$ws.Range("D41").ClearFormats()

# Row 42
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.003508'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '41CEJICEJI'

# Row 43
$ws.Range("B43").Value = 'BKEXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1038'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '42BKEXTokenBKK'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008544'
$ws.Range("D44").ClearFormats()

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005652'
$ws.Range("D45").ClearFormats()

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6815'
$ws.Range("D47").ClearFormats()

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.01588'
$ws.Range("D48").ClearFormats()

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002105'
$ws.Range("D49").ClearFormats()

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.01012'
$ws.Range("D50").ClearFormats()

